$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '25.945.66'
$dCell.Style = "Normal"
$ws.Range("E2").Value = '  -0.82%  '
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '1.643.58'
$dCell.Style = "Normal"
$ws.Range("E3").Value = '  -1.55%  '
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = '1.001'
$dCell.Style = "Normal"
$ws.Range("E4").Value = '  -0.13%  '
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '214.80'
$dCell.Style = "Normal"
$ws.Range("E5").Value = '  +2.23%  '
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '0.5209'
$dCell.Style = "Normal"
$ws.Range("E6").Value = '  -0.50%  '
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = '1.002'
$dCell.Style = "Normal"
$ws.Range("E7").Value = '  -0.06%  '
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '0.2603'
$dCell.Style = "Normal"
$ws.Range("E8").Value = '  -0.37%  '
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.06351'
$dCell.Style = "Normal"
$ws.Range("E9").Value = '  -0.04%  '
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '20.64'
$dCell.Style = "Normal"
$ws.Range("E10").Value = '  -2.23%  '
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '0.07672'
$dCell.Style = "Normal"
$ws.Range("E11").Value = '  +1.75%  '
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '1.640.57'
$dCell.Style = "Normal"
$ws.Range("E12").Value = '  -1.77%  '
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '4.421'
$dCell.Style = "Normal"
$ws.Range("E13").Value = '  -0.33%  '
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '1.867.68'
$dCell.Style = "Normal"
$ws.Range("E14").Value = '  -1.42%  '
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '0.5489'
$dCell.Style = "Normal"
$ws.Range("E15").Value = '  +0.89%  '
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₅8186'
$dCell.Style = "Normal"
$ws.Range("E16").Value = '  +2.06%  '
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '64.51'
$dCell.Style = "Normal"
$ws.Range("E17").Value = '  -2.93%  '
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '25.957.35'
$dCell.Style = "Normal"
$ws.Range("E18").Value = '  -0.91%  '
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '1.002'
$dCell.Style = "Normal"
$ws.Range("E19").Value = '  -0.11%  '
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '4.694'
$dCell.Style = "Normal"
$ws.Range("E20").Value = '  -1.08%  '
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '189.25'
$dCell.Style = "Normal"
$ws.Range("E21").Value = '  +0.82%  '
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '10.13'
$dCell.Style = "Normal"
$ws.Range("E22").Value = '  -1.54%  '
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '6.246'
$dCell.Style = "Normal"
$ws.Range("E23").Value = '  -0.24%  '
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '1.002'
$dCell.Style = "Normal"
$ws.Range("E24").Value = '  -0.09%  '
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '143.61'
$dCell.Style = "Normal"
$ws.Range("E25").Value = '  -4.05%  '
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '0.1242'
$dCell.Style = "Normal"
$ws.Range("E26").Value = '  +0.79%  '
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = '7.373'
$dCell.Style = "Normal"
$ws.Range("E27").Value = '  -1.04%  '
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '15.91'
$dCell.Style = "Normal"
$ws.Range("E28").Value = '  +0.87%  '
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '1.397'
$dCell.Style = "Normal"
$ws.Range("E29").Value = '  +2.50%  '
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '0.05880'
$dCell.Style = "Normal"
$ws.Range("E30").Value = '  -6.13%  '
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '1.262'
$dCell.Style = "Normal"
$ws.Range("E31").Value = '  -1.00%  '
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '3.392'
$dCell.Style = "Normal"
$ws.Range("E32").Value = '  -0.77%  '
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '3.388'
$dCell.Style = "Normal"
$ws.Range("E33").Value = '  -3.08%  '
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '1.641'
$dCell.Style = "Normal"
$ws.Range("E34").Value = '  -0.21%  '
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9888'
$dCell.Style = "Normal"
$ws.Range("E35").Value = '  -1.27%  '
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '2.396'
$dCell.Style = "Normal"
$ws.Range("E36").Value = '  +0.16%  '
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '2.749'
$dCell.Style = "Normal"
$ws.Range("E37").Value = '  -0.39%  '
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '0.5622'
$dCell.Style = "Normal"
$ws.Range("E38").Value = '  -5.79%  '
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '0.01598'
$dCell.Style = "Normal"
$ws.Range("E39").Value = '  -0.49%  '
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '5.852'
$dCell.Style = "Normal"
$ws.Range("E40").Value = '  -3.50%  '
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '0.8528'
$dCell.Style = "Normal"
$ws.Range("E41").Value = '  -0.66%  '
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '1.002'
$dCell.Style = "Normal"
$ws.Range("E42").Value = '  -0.20%  '
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '1.026.83'
$dCell.Style = "Normal"
$ws.Range("E43").Value = '  -7.63%  '
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '98.56'
$dCell.Style = "Normal"
$ws.Range("E44").Value = '  -2.22%  '
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '1.792.26'
$dCell.Style = "Normal"
$ws.Range("E45").Value = '  -1.46%  '
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '0.0₈107'
$dCell.Style = "Normal"
$ws.Range("E46").Value = '  -2.43%  '
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '55.50'
$dCell.Style = "Normal"
$ws.Range("E47").Value = '  -0.12%  '
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '0.9980'
$dCell.Style = "Normal"
$ws.Range("E48").Value = '  -0.48%  '
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '8.055'
$dCell.Style = "Normal"
$ws.Range("E49").Value = '  -0.17%  '
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '0.05146'
$dCell.Style = "Normal"
$ws.Range("E50").Value = '  -1.85%  '
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '0.4215'
$dCell.Style = "Normal"
$ws.Range("E51").Value = '  -0.52%  '
